$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new username value in the next empty row (row 5, since rows 1-4 are used)
$ws.Range("A5").Value = "leomessi"

# Update the selection to match the diff (D6)
$ws.Range("D6").Select()
